$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(61, "Thu Oct 17 00:27:00 2024", 128, 139),
    @(62, "Thu Oct 17 00:27:00 2024", 128, 140),
    @(63, "Thu Oct 17 00:27:01 2024", 128, 140),
    @(64, "Thu Oct 17 00:27:11 2024", 127, 138),
    @(65, "Thu Oct 17 00:27:21 2024", 127, 138),
    @(66, "Thu Oct 17 00:27:31 2024", 127, 138),
    @(67, "Thu Oct 17 00:27:41 2024", 126, 137),
    @(68, "Thu Oct 17 00:27:51 2024", 126, 137),
    @(69, "Thu Oct 17 00:28:01 2024", 126, 136),
    @(70, "Thu Oct 17 00:28:11 2024", 126, 136),
    @(71, "Thu Oct 17 00:37:01 2024", 147, 140),
    @(72, "Thu Oct 17 00:37:01 2024", 149, 142),
    @(73, "Thu Oct 17 00:37:01 2024", 153, 146),
    @(74, "Thu Oct 17 00:37:09 2024", 150, 146)
)

$startRow = 62
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $entry = $data[$i]
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
}
